# Updates the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# symbol list with refreshed values, as produced by the Sun Jan 8
# 2023 19:45 UTC GitHub Actions data-refresh run.
#
# Source cells are plain text strings (e.g. "265.79", "1.61%"), not
# numbers/percentages, so a plain Range.Value assignment would let
# Excel's autodetection reinterpret them (turning "2.14%" into a real
# percentage number, or losing exact text like "0.8500" -> 0.85).
# Force the target cell to Text first, assign, then restore the
# original "Normal" style/General format so the on-disk cell keeps
# the same look-and-feel as before, just with new text content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "267.15"
Set-TextValue $ws.Range("E2") "2.14%"
Set-TextValue $ws.Range("D3") "26.71"
Set-TextValue $ws.Range("E3") "-1.80%"
Set-TextValue $ws.Range("D4") "4.689"
Set-TextValue $ws.Range("E4") "-0.25%"
Set-TextValue $ws.Range("E5") "-1.80%"
Set-TextValue $ws.Range("D6") "6.732"
Set-TextValue $ws.Range("E6") "0.23%"
Set-TextValue $ws.Range("D7") "0.8500"
Set-TextValue $ws.Range("E7") "0.01%"
Set-TextValue $ws.Range("D8") "0.9067"
Set-TextValue $ws.Range("E8") "-1.38%"
Set-TextValue $ws.Range("D9") "0.1411"
Set-TextValue $ws.Range("E9") "-0.01%"
Set-TextValue $ws.Range("D10") "0.04993"
Set-TextValue $ws.Range("E10") "10.15%"
Set-TextValue $ws.Range("D11") "0.07086"
Set-TextValue $ws.Range("E11") "0.00%"
Set-TextValue $ws.Range("D12") "0.03176"
Set-TextValue $ws.Range("E12") "0.78%"
Set-TextValue $ws.Range("D13") "0.09020"
Set-TextValue $ws.Range("E13") "-0.44%"
Set-TextValue $ws.Range("D14") "0.001529"
Set-TextValue $ws.Range("E14") "-0.91%"
Set-TextValue $ws.Range("D15") "0.0006057"
Set-TextValue $ws.Range("E15") "-1.65%"
Set-TextValue $ws.Range("D16") "0.005974"
Set-TextValue $ws.Range("E16") "-1.01%"
Set-TextValue $ws.Range("E17") "0.06%"
Set-TextValue $ws.Range("D18") "3.169"
Set-TextValue $ws.Range("E18") "0.12%"
Set-TextValue $ws.Range("D19") "2.283"
Set-TextValue $ws.Range("E19") "4.10%"
Set-TextValue $ws.Range("D22") "4.074"
Set-TextValue $ws.Range("E22") "-0.60%"
Set-TextValue $ws.Range("D23") "0.04233"
Set-TextValue $ws.Range("E23") "-0.24%"
Set-TextValue $ws.Range("E24") "-2.42%"
Set-TextValue $ws.Range("D25") "0.004134"
Set-TextValue $ws.Range("E25") "8.74%"
Set-TextValue $ws.Range("D26") "0.0001200"
Set-TextValue $ws.Range("E26") "0.02%"
Set-TextValue $ws.Range("D27") "0.0001681"
Set-TextValue $ws.Range("E27") "5.02%"
Set-TextValue $ws.Range("D40") "0.03914"
Set-TextValue $ws.Range("E40") "-0.29%"
Set-TextValue $ws.Range("E41") "-0.04%"
Set-TextValue $ws.Range("D42") "0.004187"
Set-TextValue $ws.Range("E42") "1.32%"
Set-TextValue $ws.Range("E43") "-1.83%"
Set-TextValue $ws.Range("D44") "0.01268"
Set-TextValue $ws.Range("E44") "-8.51%"
Set-TextValue $ws.Range("D45") "0.00005135"
Set-TextValue $ws.Range("E45") "-0.59%"
Set-TextValue $ws.Range("E46") "0.01%"
Set-TextValue $ws.Range("E47") "-31.80%"
Set-TextValue $ws.Range("D48") "0.1378"
Set-TextValue $ws.Range("E48") "-17.34%"
Set-TextValue $ws.Range("E49") "0.01%"
Set-TextValue $ws.Range("E50") "0.01%"
